# Generate Report for Handoff
# Rotates the two tracked files (and their handoff/handback records) to a
# new pair of GUIDs, moves the file from "handed back" into a fresh
# "ready for handoff" cycle, and drops the now-stale handback columns
# (Latest Target File / Latest Handback File) from the per-language sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "1ba50eee-edec-4725-a32a-1788f6f95a17"
$oldGuid2 = "ae6bee11-e545-4631-b6b4-eb86fee29f82"
$oldHash1 = "5063eaa8caf0ebb66bfb39f8feabf11d67970bd2"
$oldHash2 = "0e62a27fb5cb90ba0674157b43917c4b919f89a4"

$newGuid1 = "0e71f548-af60-45c4-b839-55682e57a859"
$newGuid2 = "ffff346c3780-2602-4cf0-9d3f-bc4d5f73ec6c"
$newHash  = "0aa05150851e84d6a81556a899c19be5e8830348"

$status = "Ready for handoff"
$overviewDate = "2016-46-11 09:46:47"
$zhHandoffDatetime = "2016-03-11 09:46:44"
$deHandoffDatetime = "2016-03-11 09:46:47"
$staleHandbackDate = "0001-01-01 00:00:00"

$md1 = "$newGuid1.md"
$md2 = "$newGuid2.md"
$zhXlf = "$newGuid1.$newHash.zh-cn.xlf"
$deXlf = "$newGuid1.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md1", "", "", $md1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md2", "", "", $md2)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $md1
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $status
$wsZh.Range("D2").Value = $zhXlf
$wsZh.Range("E2").Value = $zhHandoffDatetime
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $staleHandbackDate
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $md2
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $zhXlf
$wsZh.Range("E3").Value = $zhHandoffDatetime
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $staleHandbackDate
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md1", "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md1", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6a41b1556de56efb39c5a2d4bdb31cc436b6acf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md2", "", "", $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md2", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6a41b1556de56efb39c5a2d4bdb31cc436b6acf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $md1
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $status
$wsDe.Range("D2").Value = $deXlf
$wsDe.Range("E2").Value = $deHandoffDatetime
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $staleHandbackDate
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $md2
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $deXlf
$wsDe.Range("E3").Value = $deHandoffDatetime
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $staleHandbackDate
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md1", "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md1", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/add6da2ed3c3392a1980e52571b17d0b0e8050b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md2", "", "", $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/30cec44aa0f4f4230cc051b664dbcb49639543aa/e2e/$md2", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/add6da2ed3c3392a1980e52571b17d0b0e8050b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)
